$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 06:07"

# --- Pakistan (row 20) ---
$ws.Range("B20").Value = 302020
$ws.Range("C20").Value = 539
$ws.Range("D20").Value = 289806
$ws.Range("E20").Value = 5831
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 6383

# --- Honduras (row 50) ---
$ws.Range("B50").Value = 67789
$ws.Range("C50").Value = 653
$ws.Range("D50").Value = 17873
$ws.Range("E50").Value = 47837
$ws.Range("G50").Value = 14
$ws.Range("H50").Value = 2079

# --- Tailandia (row 128) ---
$ws.Range("B128").Value = 3475
$ws.Range("C128").Value = 2
$ws.Range("E128").Value = 105

# --- Belice / Letonia swap positions (rows 159-160) because Belice's
#     totals overtook Letonia's after the update, re-sorting the table ---
# Row 159 keeps the "Belice" label and gets the freshly updated figures.
$ws.Range("A159").Value = "Belice"
$ws.Range("B159").Value = 1480
$ws.Range("C159").Value = 22
$ws.Range("D159").Value = 492
$ws.Range("E159").Value = 969
$ws.Range("H159").Value = 19

# Row 160 now carries the "Letonia" label with its previous (unchanged) figures.
$ws.Range("A160").Value = "Letonia"
$ws.Range("B160").Value = 1474
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 1248
$ws.Range("E160").Value = 191
$ws.Range("H160").Value = 35

# --- Mongolia (row 184) ---
$ws.Range("D184").Value = 300
$ws.Range("E184").Value = 11
